$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A so the existing Category/Description/Amount
# columns shift right by one (A->B, B->C, C->D), making room for a new
# date column in A.
$ws.Range("A1:A2").EntireColumn.Insert()

# Populate the new date column (stored as plain text, like the original
# writer produced) for both movement rows.
$ws.Range("A1").Value2 = "27/04/2018"
$ws.Range("A2").Value2 = "27/04/2018"
